$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D (shifts existing D:K data to E:L for the newest reporting period)
$ws.Columns("D").Insert()

# Copy cell formatting (number format/style) from the adjacent column E into the
# newly inserted column D, scoped to the three data blocks on the sheet
# (Income Statement, Balance Sheet, Cash Flow Statement) so that section-header
# rows (5, 6, 37, 79) that have no data cells stay untouched.
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)
$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)
$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new column D with the latest reporting period's figures.
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 11221000
$ws.Range("D9").Value = 6873000
$ws.Range("D10").Value = 4348000
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 33000
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 10056000
$ws.Range("D18").Value = 1165000
$ws.Range("D20").Value = 4000
$ws.Range("D21").Value = 1426000
$ws.Range("D22").Value = 88000
$ws.Range("D23").Value = 1081000
$ws.Range("D24").Value = 258000
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 823000
$ws.Range("D27").Value = 776000
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -4000
$ws.Range("D33").Value = 776000
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 776000
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 538000
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 1385000
$ws.Range("D44").Value = 1541000
$ws.Range("D45").Value = 93000
$ws.Range("D46").Value = 3557000
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 1352000
$ws.Range("D49").Value = 952000
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 12000
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 5873000
$ws.Range("D57").Value = 678000
$ws.Range("D58").Value = 130000
$ws.Range("D59").Value = 693000
$ws.Range("D60").Value = 1501000
$ws.Range("D61").Value = 2090000
$ws.Range("D62").Value = 189000
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 3952000
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 7869000
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 1921000
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 776000
$ws.Range("D83").Value = 257000
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 1057000
$ws.Range("D91").Value = -239000
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -166000
$ws.Range("D96").Value = -316000
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -670000
$ws.Range("D101").Value = -10000
$ws.Range("D102").Value = 211000
